$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
# Column Q ("weaponSlot") is populated with 16 for every data row (5-35).
# Row 17 previously had a stray Q/R pair (Q17=10, R17="10|5|6"); Q17 becomes
# 16 like every other row and R17 is cleared out entirely.
for ($r = 5; $r -le 35; $r++) {
    $ws.Range("Q$r").Value = 16
}
$ws.Range("R17").ClearContents()

# --- View / selection (best effort) --------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A4").Select()
$ws.Range("R36").Select()
